$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "2025-04-28 22:30:35"
$ws.Range("B26").Value = 476
